$wb = $excel.ActiveWorkbook

# Sheet "展览" (rId1 / sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2966
$ws1.Range("F3").Value = 6416
$ws1.Range("F12").Value = 7446
$ws1.Range("F19").Value = 9098
$ws1.Range("F39").Value = 3897
$ws1.Range("F45").Value = 239
$ws1.Range("F48").Value = 29

# Sheet "全部类型" (rId4 / sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2966
$ws4.Range("F6").Value = 6416
$ws4.Range("F18").Value = 7446
$ws4.Range("F24").Value = 9098
$ws4.Range("F41").Value = 3897
$ws4.Range("F46").Value = 239
$ws4.Range("F48").Value = 29
